$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new reporting-period column before column D.
# This shifts the existing D:K data one column to the right (-> E:L),
# matching the widened data range (A5:K102 -> A5:L102).
$ws.Columns("D:D").Insert()

# The freshly inserted column D is blank and carries column C formatting;
# re-apply the correct number formats/styles by copying them from column E
# (which now holds what used to be column D).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Give column D the same width as the rest of the data columns.
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# --- Populate the new column D with the latest reporting-period figures.
# A handful of subtotal/total rows were also restated for the two most
# recent prior periods, so columns E and F are corrected there too.

$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 300200
$ws.Range("E8").Value = 154600
$ws.Range("F8").Value = 66100
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 2500
$ws.Range("D15").Value = 2600
$ws.Range("D17").Value = 39400
$ws.Range("E17").Value = 9700
$ws.Range("F17").Value = 3200
$ws.Range("D18").Value = 260800
$ws.Range("E18").Value = 144900
$ws.Range("F18").Value = 62900
$ws.Range("D20").Value = 900
$ws.Range("D21").Value = 264300
$ws.Range("E21").Value = 145400
$ws.Range("F21").Value = 62900
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 261800
$ws.Range("E23").Value = 145100
$ws.Range("F23").Value = 62900
$ws.Range("D24").Value = 52000
$ws.Range("E24").Value = 47900
$ws.Range("F24").Value = 20600
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 209700
$ws.Range("E26").Value = 97100
$ws.Range("F26").Value = 42300
$ws.Range("D27").Value = 209700
$ws.Range("E27").Value = 97100
$ws.Range("F27").Value = 42300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = 100
$ws.Range("F29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -900
$ws.Range("D33").Value = 209700
$ws.Range("E33").Value = 97200
$ws.Range("F33").Value = 42300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 209700
$ws.Range("E35").Value = 97200
$ws.Range("F35").Value = 42300
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 119600
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 48800
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 9400
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 99600
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = "NA"
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 285100
$ws.Range("D57").Value = 10500
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 1600
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 14900
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 40400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 244700
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 209700
$ws.Range("E81").Value = 97200
$ws.Range("F81").Value = 42300
$ws.Range("D83").Value = 2600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 191600
$ws.Range("D91").Value = -81600
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -81500
$ws.Range("D96").Value = -31700
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -70000
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 40100
